# Insert a new weekly record as the first entry (row 45) of the Ají
# (Hortaliza) dataset for Vega Monumental Concepción, pushing the existing
# rows 45-99 down to 46-100.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; everything below (old rows 45-99) shifts
# down to 46-100, carrying its formatting (incl. the date style on column D).
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the new weekly record. Columns that are
# identical to the row that used to occupy this slot (A, B, C, E, F, G, I,
# N, Q, R) are simply re-entered with the same values.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44679
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112021
$ws.Range("G45").Value = "Ají"
$ws.Range("H45").Value = "Chilena(o)"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 20000
$ws.Range("L45").Value = 22000
$ws.Range("M45").Value = 21000
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Región Metropolitana"
$ws.Range("P45").Value = 840
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
